# Adds the new daily COVID-19 data row (2020-06-10, serial 43992) to the
# bottom of the "Covid-19 podatki" sheet / "Tabela1" table, matching the
# upstream data-bot commit that appended row 92.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Grow the table by one row so the table ref / autoFilter ref and the
#    sheet dimension all extend from ...J91 to ...J92, same as Excel does
#    when a new row is appended right under an existing ListObject.
$lo = $ws.ListObjects.Item(1)
[void]$lo.ListRows.Add()

# 2) Row 72 is an existing "highlighted" data row (alternating banding +
#    border), which is the same formatting the new last row should carry.
#    Copy just its formatting down onto the freshly-added row 92 so the
#    cells pick up the existing style entries instead of Excel inventing
#    new ones.
$ws.Range("A72:J72").Copy()
$ws.Range("A92:J92").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Write the new day's figures.
$ws.Range("A92").Value = 43992
$ws.Range("B92").Value = 85626
$ws.Range("C92").Value = 758
$ws.Range("D92").Value = 1488
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 6
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 109
$ws.Range("J92").Value = 0

# 4) Move the selection to the newly entered row, matching the saved
#    workbook's cursor position.
[void]$ws.Range("A92:J92").Select()
